$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Key", "Meaning"),
    @("NN", "Neural Networks"),
    @("G", "Git"),
    @("AWS", "Amazon Web Services"),
    @("DA", "Data Analysis"),
    @("AZ", "Azure"),
    @("DB", "Databases"),
    @("AF", "Airflow"),
    @("KB", "Kubernetes"),
    @("V", "Data Visualization"),
    @("FT", "FinTech"),
    @("P", "Python"),
    @("R", "R"),
    @("J", "Java"),
    @("JS", "JavaScript"),
    @("C", "C/C++"),
    @("ML", "Machine Learning"),
    @("NL", "Natural Language Processing/Understanding"),
    @("ST", "Statistics"),
    @("M", "Mathematics"),
    @("MA", "Matlab"),
    @("DW", "Data Wrangling"),
    @("SQ", "SQL"),
    @("LD", "Large Datasets"),
    @("XC", "Excel"),
    @("TB", "Tableau"),
    @("SE", "Software Engineering"),
    @("SP", "Spark"),
    @("HV", "Hive"),
    @("HP", "Hadoop"),
    @("CF", "Caffe2"),
    @("NG", "Networks/Graphs"),
    @("NP", "Numpy"),
    @("SC", "Scipy"),
    @("PD", "Pandas"),
    @("SK", "SKLearn"),
    @("TF", "TensorFlow"),
    @("PT", "PyTorch"),
    @("TH", "Theano"),
    @("K", "Keras"),
    @("TS", "Time Series"),
    @("OO", "Object-oriented programming"),
    @("DL", "Deep Learning"),
    @("PB", "Publications"),
    @("BQ", "BigQuery"),
    @("DE", "Data Engineering"),
    @("SL", "Scala"),
    @("UL", "Unix/Linux"),
    @("CA", "Cassandra"),
    @("CS", "Computer Science"),
    @("D3", "D3"),
    @("DK", "Docker"),
    @("HB", "Hbase"),
    @("MG", "MongoDB"),
    @("NS", "NoSQL"),
    @("PE", "Perl"),
    @("PG", "Pig"),
    @("SAS", "SAS"),
    @("SS", "SPSS"),
    @("SD", "Software Development"),
    @("PM", "Project Management")
)

$rowCount = $data.Count
$colCount = 2

$values = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $values[$i, 0] = $data[$i][0]
    $values[$i, 1] = $data[$i][1]
}

# A1 is temporarily left as its original "Key" value (row 0 of $data) so that the bulk
# write below does not register the brand-new "Code" string ahead of the other new
# strings introduced further down the sheet (rows 49-61). We overwrite A1 afterwards.
$range = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($rowCount, 2))
$range.Value = $values

# Now set A1 to its final value last, so "Code" is appended to the shared
# strings table after all the other newly-introduced strings.
$ws.Cells.Item(1, 1).Value = "Code"

$ws.Range("A2").Select()

Write-Host "Wrote $rowCount rows"
